# Updated cryptos list on Thu Jan 11 19:23:50 UTC 2024 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) columns for every coin row with
# the latest pull from the data source.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep numeric-looking price strings stored as text (matches the source feed,
# which writes plain text cells), so Excel does not silently reinterpret them
# as numbers when the new price happens to parse as a float.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the refreshed values
$ws.Range("D2").Value = '46.423.24'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '2.603.46'
$ws.Range("E3").Value = '  +5.76%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '307.59'
$ws.Range("E5").Value = '  +3.28%  '
$ws.Range("D6").Value = '100.27'
$ws.Range("E6").Value = '  +2.86%  '
$ws.Range("D7").Value = '0.603'
$ws.Range("E7").Value = '  +4.21%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +12.07%  '
$ws.Range("E10").Value = '  +10.77%  '
$ws.Range("D11").Value = '54.70'
$ws.Range("E11").Value = '  +1.85%  '
$ws.Range("D12").Value = '0.0846'
$ws.Range("E12").Value = '  +7.19%  '
$ws.Range("D13").Value = '8.18'
$ws.Range("E13").Value = '  +13.50%  '
$ws.Range("D14").Value = '3.004.50'
$ws.Range("E14").Value = '  +5.90%  '
$ws.Range("E15").Value = '  +1.41%  '
$ws.Range("D16").Value = '2.616.05'
$ws.Range("E16").Value = '  +6.62%  '
$ws.Range("D17").Value = '0.922'
$ws.Range("E17").Value = '  +7.33%  '
$ws.Range("D18").Value = '14.99'
$ws.Range("E18").Value = '  +5.74%  '
$ws.Range("D19").Value = '46.579.73'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  +7.06%  '
$ws.Range("D21").Value = '12.99'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").Value = '6.78'
$ws.Range("E22").Value = '  +7.62%  '
$ws.Range("D23").Value = '279.98'
$ws.Range("E23").Value = '  +13.43%  '
$ws.Range("D24").Value = '71.78'
$ws.Range("E24").Value = '  +5.83%  '
$ws.Range("D25").Value = '3.04'
$ws.Range("E25").Value = '  +7.69%  '
$ws.Range("E26").Value = '  +9.41%  '
$ws.Range("D27").Value = '28.75'
$ws.Range("E27").Value = '  +33.35%  '
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("D29").Value = '4.03'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").Value = '10.65'
$ws.Range("E30").Value = '  +7.68%  '
$ws.Range("D31").Value = '39.19'
$ws.Range("E31").Value = '  -2.82%  '
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").Value = '6.38'
$ws.Range("E33").Value = '  +13.21%  '
$ws.Range("D34").Value = '3.66'
$ws.Range("E34").Value = '  -6.29%  '
$ws.Range("D35").Value = '2.83'
$ws.Range("E35").Value = '  +2.93%  '
$ws.Range("E36").Value = '  +8.17%  '
$ws.Range("E37").Value = '  +7.54%  '
$ws.Range("D38").Value = '151.80'
$ws.Range("E38").Value = '  +2.06%  '
$ws.Range("E39").Value = '  +7.71%  '
$ws.Range("E40").Value = '  +5.39%  '
$ws.Range("D41").Value = '23.07'
$ws.Range("E41").Value = '  +38.57%  '
$ws.Range("D42").Value = '16.17'
$ws.Range("E42").Value = '  +4.53%  '
$ws.Range("D43").Value = '0.0332'
$ws.Range("E43").Value = '  +9.05%  '
$ws.Range("D44").Value = '3.64'
$ws.Range("E44").Value = '  +9.20%  '
$ws.Range("E45").Value = '  +2.38%  '
$ws.Range("D46").Value = '2.131.79'
$ws.Range("E46").Value = '  +7.13%  '
$ws.Range("D47").Value = '0.999'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").Value = '93.20'
$ws.Range("E48").Value = '  +0.57%  '
$ws.Range("D49").Value = '9.56'
$ws.Range("E49").Value = '  +10.53%  '
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("D51").Value = '109.66'
$ws.Range("E51").Value = '  +7.00%  '
